$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Expected Life(pre)" column (column D), which always held 0.
# This shifts "Expected Life(post)" into column D and "Depriciation Method" into column E.
$ws.Range("D1").EntireColumn.Delete()

# Rename the (now) column D header from "Expected Life(post)" to "Expected Life"
$ws.Range("D1").Value = "Expected Life"

# Update the Expected Life values per the corrected asset schedule
$ws.Range("D2").Value = 10
$ws.Range("D3").Value = 5
$ws.Range("D4").Value = 10
$ws.Range("D5").Value = 10
$ws.Range("D6").Value = 10
$ws.Range("D7").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("D9").Value = 10
$ws.Range("D10").Value = 15
$ws.Range("D11").Value = 5
$ws.Range("D12").Value = 10
$ws.Range("D13").Value = 10
$ws.Range("D14").Value = 5

# Update the selection to match the new sheet state
$null = $ws.Range("I5:J5").Select()
